# Update "想去人数" (interested count) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 39
$ws1.Range("F10").Value = 139
$ws1.Range("F11").Value = 10193
$ws1.Range("F13").Value = 267
$ws1.Range("F16").Value = 11806
$ws1.Range("F17").Value = 12212

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 39
$ws4.Range("F11").Value = 139
$ws4.Range("F12").Value = 10193
$ws4.Range("F14").Value = 267
$ws4.Range("F17").Value = 11806
$ws4.Range("F18").Value = 12212
